# Update countries & provincias Spain
# Updates the "Datos actualizados" timestamp and refreshes the COVID case
# counters for several countries. "Costa de Marfil" overtakes
# "Consejo Danes para los Refugiados" in total cases, so the two rows swap
# their country labels (the table is kept sorted by total cases descending)
# while "Consejo Danes para los Refugiados" keeps its (unchanged) prior
# figures one row further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 22:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1739609
$ws.Range("C4").Value = 14334
$ws.Range("E4").Value = 1154801
$ws.Range("G4").Value = 1144
$ws.Range("H4").Value = 101716

# Alemania (row 11)
$ws.Range("B11").Value = 181821
$ws.Range("C11").Value = 533
$ws.Range("E11").Value = 10488
$ws.Range("G11").Value = 35
$ws.Range("H11").Value = 8533

# Row 83: now "Costa de Marfil" (overtakes Consejo Danes in total cases)
$ws.Range("A83").Value = "Costa de Marfil"
$ws.Range("B83").Value = 2556
$ws.Range("C83").Value = 79
$ws.Range("D83").Value = 1302
$ws.Range("E83").Value = 1223
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 31

# Row 84: now "Consejo Danes para los Refugiados" (unchanged totals, drops a row)
$ws.Range("A84").Value = "Consejo Danes para los Refugiados"
$ws.Range("B84").Value = 2546
$ws.Range("C84").Value = 143
$ws.Range("D84").Value = 365
$ws.Range("E84").Value = 2113
$ws.Range("H84").Value = 68

# Republica de Africa Central (row 131)
$ws.Range("B131").Value = 702
$ws.Range("C131").Value = 31
$ws.Range("D131").Value = 23
$ws.Range("E131").Value = 678

# Santo Tome y Principe (row 139)
$ws.Range("B139").Value = 443
$ws.Range("C139").Value = 2
$ws.Range("E139").Value = 363
